# Applies the Homework-4 / Quiz-4 formatting update:
#   1. Splits the "You might need to use a for loop..." run in two and
#      drops a fresh _GoBack bookmark at the split point (Word keeps only
#      one _GoBack bookmark, so adding this one also removes the old one
#      that currently sits after the screenshot).
#   2. Resizes the inline screenshot picture (both the display extent and
#      the underlying picture-frame extent).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: resize the inline picture.
# We rebuild its <w:drawing> with distinct wp:extent / a:ext values (the
# plain Width/Height setters always keep those two in lock-step, but the
# authored diff wants them slightly different), so splice in literal XML
# for just that drawing run.
# ---------------------------------------------------------------------
$shp = $d.InlineShapes.Item(1)
$insertPos = $shp.Range.Start
$shp.Delete()

$drawingXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:noProof/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="433D211F" wp14:editId="1B2C79B5"><wp:extent cx="4421951" cy="701040"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name="Screen Shot 2018-12-09 at 11.38.03 AM.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="4522541" cy="716987"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Range($insertPos, $insertPos)
$target.InsertXML($drawingXml)

# ---------------------------------------------------------------------
# Step 2: split the "for loop" sentence run and drop the _GoBack bookmark
# right after "...print out th". Adding a bookmark named _GoBack moves
# Word's single hidden _GoBack bookmark here, removing it from its old
# spot after the picture automatically.
# ---------------------------------------------------------------------
$sentence = "You might need to use a " + [char]0x201C + "for" + [char]0x201D + " loop to print out the items from your dictionary. Here is one way you can do that:"
$findRange = $d.Content
$findRange.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$splitOffset = $findRange.Start + 50
$splitPoint = $d.Range($splitOffset, $splitOffset)
$splitPoint.Bookmarks.Add("_GoBack") | Out-Null
